$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reference cells that already carry the "ship" (style 6) and
# --- "hit marker" (style 9) formatting used elsewhere on the board.
$shipFormatSrc = $ws.Range("T11")   # style 6: red fill, thin border
$hitFormatSrc  = $ws.Range("U12")   # style 9: bold font, red fill, thin border

# Row 9: R9:T9 become ship cells (value 3, style 6)
$ws.Range("R9").Value = 3
$ws.Range("S9").Value = 3
$ws.Range("T9").Value = 3
$shipFormatSrc.Copy()
$ws.Range("R9:T9").PasteSpecial(-4122)

# Row 10: R10 & T10 ship cells, S10 becomes a hit marker ("X")
$ws.Range("R10").Value = 3
$ws.Range("T10").Value = 3
$shipFormatSrc.Copy()
$ws.Range("R10").PasteSpecial(-4122)
$ws.Range("T10").PasteSpecial(-4122)

$ws.Range("S10").Value = "X"
$hitFormatSrc.Copy()
$ws.Range("S10").PasteSpecial(-4122)

# Row 11: R11 becomes a ship cell, S11 becomes a hit marker ("X")
$ws.Range("R11").Value = 3
$shipFormatSrc.Copy()
$ws.Range("R11").PasteSpecial(-4122)

$ws.Range("S11").Value = "X"
$hitFormatSrc.Copy()
$ws.Range("S11").PasteSpecial(-4122)

# Row 12: R12 & S12 become ship cells
$ws.Range("R12").Value = 3
$ws.Range("S12").Value = 3
$shipFormatSrc.Copy()
$ws.Range("R12:S12").PasteSpecial(-4122)

# Row 14: U14 changes from a ship cell into a hit marker ("X")
$ws.Range("U14").Value = "X"
$hitFormatSrc.Copy()
$ws.Range("U14").PasteSpecial(-4122)

# Row 15/16: a second ship was placed overlapping the first one - these
# cells get the new "overlap" fill (same red pattern, no white background).
# Build the two new styles cleanly on a single cell each, then fan the
# resulting formatting out to the rest of the affected cells.
$ws.Range("T15").Value = 3
$ws.Range("T15").Interior.Pattern = 17
$ws.Range("T15").Interior.Color = 255

$ws.Range("U15").Value = "X"
$ws.Range("U15").Font.Bold = $true
$ws.Range("U15").Interior.Pattern = 17
$ws.Range("U15").Interior.Color = 255

$ws.Range("V15").Value = 3
$ws.Range("T15").Copy()
$ws.Range("V15").PasteSpecial(-4122)

$ws.Range("T16").Value = 3
$ws.Range("U16").Value = 3
$ws.Range("V16").Value = 3
$ws.Range("T15").Copy()
$ws.Range("T16:V16").PasteSpecial(-4122)

# ToDo list: the finished item's text is replaced by the next one, and
# a new open item is appended at the end of the list.
$ws.Range("P23").Value = "'- beim drehen auf andere Schiffe konrollieren"
$ws.Range("P29").Value = "'- Gesperrter Bereich mein Setzen schreibt über grenze"

# Restore the active selection as recorded in the saved workbook.
$ws.Range("U23").Select() | Out-Null
